$d = $word.ActiveDocument

# 1. Merge "Ik vond dat ik mijn taken..." run fragments (user stories / wordpress)
#    back into a single run, dropping the spell-check proofErr markers that
#    wrapped "stories" and "wordpress".
$d.Content.Find.Execute(
    "Ik vond dat ik mijn taken goed onder kon verdelen in user stories, wordpress had ik goed onder de knie dus ik liep niet vaak ergens heel erg tegen aan.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Ik vond dat ik mijn taken goed onder kon verdelen in user stories, wordpress had ik goed onder de knie dus ik liep niet vaak ergens heel erg tegen aan.",
    2) | Out-Null

# 2. Insert the missing reflection sentence about the responsive website,
#    right after the "Ik had mijn dagen..." paragraph.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Ik had mijn dagen net wat beter kunnen onderverdelen in taken in plaats van allemaal taken hebben en die gewoon beginnen uit te voeren.") {
        $p.Range.InsertParagraphAfter()
        $newP = $d.Paragraphs.Item($i + 1)
        $newP.Range.InsertAfter("Ik had graag de website mobiel responsive willen maken als ik hier meer tijd voor kreeg")
        break
    }
}

# 3. Merge "Goede samenwerking en duidelijke communicatie" back into one run.
$d.Content.Find.Execute(
    "Goede samenwerking en duidelijke communicatie",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Goede samenwerking en duidelijke communicatie",
    2) | Out-Null
